$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct the May 2024 value (F6) and add the June 2024 value (F7)
$ws.Range("F6").Value = 93775
$ws.Range("F7").Value = 106428

# Select F8, matching the sheetView selection recorded after data entry
$ws.Range("F8").Select()
